$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.266.92"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.369.34"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.41"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.10"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.379.17"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0989"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  +6.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.327"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.792.64"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.257.51"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.50"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.377.57"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.01"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.54"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.27"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.75"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.19"
$ws.Range("E28").Value = "  -3.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.70"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0711"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.85"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.58"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.69"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.826"
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("E41").Value = "  -3.73%  "
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "125.24"
$ws.Range("E43").Value = "  -5.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.74"
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0898"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "236.70"
$ws.Range("E47").Value = "  -6.39%  "
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0206"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.98"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.954"
$ws.Range("E51").Value = "  +0.17%  "
